$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 112, shifting existing rows 112:237 down to 113:238
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with the new record
$ws.Cells.Item(112, 1).Value = 5
$ws.Cells.Item(112, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value = "Maule"
$ws.Cells.Item(112, 4).Value = 44539
$ws.Cells.Item(112, 5).Value = 7
$ws.Cells.Item(112, 6).Value = 100114013
$ws.Cells.Item(112, 7).Value = "Zanahoria"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 350
$ws.Cells.Item(112, 11).Value = 8000
$ws.Cells.Item(112, 12).Value = 8000
$ws.Cells.Item(112, 13).Value = 8000
$ws.Cells.Item(112, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(112, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(112, 16).Value = 400
$ws.Cells.Item(112, 17).Value = 20
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same numeric date style (s="2") as the other date cells
$ws.Cells.Item(112, 4).NumberFormat = $ws.Cells.Item(113, 4).NumberFormat()
